$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 52, shifting existing rows 52:70 down to 53:71
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new weekly price record
$ws.Range("A52").Value = 10
$ws.Range("B52").Value = "Vega Modelo de Temuco"
$ws.Range("C52").Value = "La Araucanía"
$ws.Range("D52").Value = 44839
$ws.Range("E52").Value = 9
$ws.Range("F52").Value = 300000000
$ws.Range("G52").Value = "Espárragos"
$ws.Range("H52").Value = "Sin especificar"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 300
$ws.Range("K52").Value = 1700
$ws.Range("L52").Value = 1700
$ws.Range("M52").Value = 1700
$ws.Range("N52").Value = "$/kilo"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 1700
$ws.Range("Q52").Value = 1
$ws.Range("R52").Value = "Hortaliza"
